# New Test Case Added
# Replaces the old sign-in / create-account test fixtures (sap198x / abctestemailN
# addresses, QqwertyQ@123 password) with a new batch of test emails
# (abctestemail237{8,9,0,1,2}!!!!@gmail.com) and a strengthened password
# (QqwertyQ@123!). Also updates the active sheet / selection to what was left
# selected when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "signin"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("signin")
$ws1.Cells.Hyperlinks.Delete()

$ws1.Range("B2").Value = "QqwertyQ@123!"
$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:QqwertyQ@123!") | Out-Null

$ws1.Range("B3").Value = "QqwertyQ@123!"
$ws1.Hyperlinks.Add($ws1.Range("B3"), "mailto:QqwertyQ@123!") | Out-Null

$ws1.Range("A6").Value = "abctestemail2372!!!!@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("A6"), "mailto:abctestemail2372!!!!@gmail.com") | Out-Null

$ws1.Range("B4").Value = "QqwertyQ@123!"
$ws1.Hyperlinks.Add($ws1.Range("B4"), "mailto:QqwertyQ@123!") | Out-Null

$ws1.Range("B5").Value = "QqwertyQ@123!"
$ws1.Hyperlinks.Add($ws1.Range("B5"), "mailto:QqwertyQ@123!") | Out-Null

$ws1.Range("B6").Value = "QqwertyQ@123!"
$ws1.Hyperlinks.Add($ws1.Range("B6"), "mailto:QqwertyQ@123!") | Out-Null

$ws1.Range("A2:A5").Value = "abctestemail237!!@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("A2:A5"), "mailto:abctestemail237!!@gmail.com", "", "", "abctestemail237!!@gmail.com") | Out-Null

$ws1.Range("A2").Value = "abctestemail2378!!!!@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:abctestemail2378!!!!@gmail.com") | Out-Null

$ws1.Range("A3").Value = "abctestemail2379!!!!@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("A3"), "mailto:abctestemail2379!!!!@gmail.com") | Out-Null

$ws1.Range("A4").Value = "abctestemail2370!!!!@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("A4"), "mailto:abctestemail2370!!!!@gmail.com") | Out-Null

$ws1.Range("A5").Value = "abctestemail2371!!!!@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("A5"), "mailto:abctestemail2371!!!!@gmail.com") | Out-Null

$ws1.Range("A2:B6").Style = "Hyperlink"

# ----------------------------------------------------------------------
# Sheet "createaccount"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("createaccount")
$ws2.Cells.Hyperlinks.Delete()

$ws2.Range("D2").Value = "QqwertyQ@123!"
$ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:QqwertyQ@123!") | Out-Null
$ws2.Range("D2").NumberFormat = "@"

$ws2.Range("D3").Value = "QqwertyQ@123!"
$ws2.Hyperlinks.Add($ws2.Range("D3"), "mailto:QqwertyQ@123!") | Out-Null
$ws2.Range("D3").NumberFormat = "@"

$ws2.Range("D4").Value = "QqwertyQ@123!"
$ws2.Range("D5").Value = "QqwertyQ@123!"
$ws2.Range("D6").Value = "QqwertyQ@123!"

$ws2.Range("C6").Value = "abctestemail2372!!!!@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("C6"), "mailto:abctestemail2372!!!!@gmail.com") | Out-Null

$ws2.Range("C2:C5").Value = "abctestemail237!!@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("C2:C5"), "mailto:abctestemail237!!@gmail.com", "", "", "abctestemail237!!@gmail.com") | Out-Null

$ws2.Range("C2").Value = "abctestemail2378!!!!@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:abctestemail2378!!!!@gmail.com") | Out-Null

$ws2.Range("C3").Value = "abctestemail2379!!!!@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("C3"), "mailto:abctestemail2379!!!!@gmail.com") | Out-Null

$ws2.Range("C4").Value = "abctestemail2370!!!!@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("C4"), "mailto:abctestemail2370!!!!@gmail.com") | Out-Null

$ws2.Range("C5").Value = "abctestemail2371!!!!@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("C5"), "mailto:abctestemail2371!!!!@gmail.com") | Out-Null

$ws2.Range("C2:C6").Style = "Hyperlink"

# ----------------------------------------------------------------------
# Sheet "searchdata" (values unchanged - only selection moves)
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("searchdata")

# ----------------------------------------------------------------------
# Selection / active-sheet bookkeeping (match the saved view state)
# ----------------------------------------------------------------------
$ws1.Range("A6").Select()
$ws3.Range("A10").Select()
$ws2.Range("C6").Select()
$ws2.Activate()
